$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.193.27"
$ws.Range("E2").Value = "  +1.62%  "

# Row 3
$ws.Range("D3").Value = "1.643.92"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").Value = "'217.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6
$ws.Range("E6").Value = "  +1.14%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("E8").Value = "  +0.33%  "

# Row 9
$ws.Range("E9").Value = "  +0.35%  "

# Row 10
$ws.Range("D10").Value = "'19.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.17%  "

# Row 11
$ws.Range("D11").Value = "'0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "

# Row 12
$ws.Range("D12").Value = "1.872.86"
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
$ws.Range("D13").Value = "'4.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.73%  "

# Row 14
$ws.Range("D14").Value = "1.655.78"
$ws.Range("E14").Value = "  +1.23%  "

# Row 15
$ws.Range("E15").Value = "  -3.21%  "

# Row 16
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  +0.04%  "

# Row 17
$ws.Range("D17").Value = "'63.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").Value = "26.194.61"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("D19").Value = "'0.999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "

# Row 20
$ws.Range("E20").Value = "  -0.81%  "

# Row 21
$ws.Range("D21").Value = "'194.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "

# Row 22
$ws.Range("D22").Value = "'10.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("E23").Value = "  -0.61%  "

# Row 24
$ws.Range("D24").Value = "'1.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.35%  "

# Row 25
$ws.Range("E25").Value = "  -0.19%  "

# Row 26
$ws.Range("D26").Value = "'142.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.64%  "

# Row 27
$ws.Range("E27").Value = "  +0.93%  "

# Row 28
$ws.Range("D28").Value = "'6.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.50%  "

# Row 29
$ws.Range("D29").Value = "'15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.72%  "

# Row 30
$ws.Range("E30").Value = "  +0.75%  "

# Row 31
$ws.Range("E31").Value = "  +1.77%  "

# Row 32
$ws.Range("E32").Value = "  +0.86%  "

# Row 33
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("E34").Value = "  +1.53%  "

# Row 35
$ws.Range("D35").Value = "'2.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.54%  "

# Row 36
$ws.Range("E36").Value = "  +0.85%  "

# Row 37
$ws.Range("D37").Value = "1.133.28"
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("E38").Value = "  +1.49%  "

# Row 39
$ws.Range("E39").Value = "  -1.14%  "

# Row 40
$ws.Range("E40").Value = "  +1.13%  "

# Row 41
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.26%  "

# Row 42
$ws.Range("D42").Value = "'100.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "

# Row 43
$ws.Range("E43").Value = "  -1.30%  "

# Row 44
$ws.Range("D44").Value = "'0.797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "

# Row 45
$ws.Range("D45").Value = "1.781.60"

# Row 46
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -0.30%  "

# Row 47
$ws.Range("D47").Value = "'56.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.76%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.43%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.86%  "

# Row 50
$ws.Range("E50").Value = "  +0.12%  "

# Row 51
$ws.Range("E51").Value = "  +2.74%  "
